$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Profile72" test-case row (row 73) ---------------------
# Copy the formatting of the last existing row (72) down onto the new
# row so the new cells pick up the same styles (borders / wrap text).
$ws.Range("A72:E72").Copy()
$ws.Range("A73:E73").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A73").Value = "Profile72"
$ws.Range("C73").Value = "Verify that Watchlist tab infinite scroll displaying the more available records"
$ws.Range("B73").Value = "OPQA-4820"
$ws.Range("D73").Value = "Y"

# --- Fill in the JIRA id that was still a placeholder on row 72 ---------
$ws.Range("B72").Value = "OPQA-4821"

# --- Update the view state to match where the author left the cursor ----
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
$ws.Range("C65").Select() | Out-Null
